# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.718.97'
$ws.Range('E2').Value = '  +2.41%  '
$ws.Range('D3').Value = '1.890.68'
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '247.42'
$ws.Range('E5').Value = '  +1.81%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4943'
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2959'
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06819'
$ws.Range('E9').Value = '  +2.56%  '
$ws.Range('D10').Value = '1.890.74'
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('E11').Value = '  +2.52%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '92.13'
$ws.Range('E12').Value = '  +6.36%  '
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6792'
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.068'
$ws.Range('E15').Value = '  +2.92%  '
$ws.Range('D16').Value = '30.670.40'
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000007973'
$ws.Range('E17').Value = '  +1.24%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.26'
$ws.Range('E19').Value = '  +3.51%  '
$ws.Range('D20').Value = '2.135.01'
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.841'
$ws.Range('E22').Value = '  +1.31%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '190.76'
$ws.Range('E23').Value = '  +33.59%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.068'
$ws.Range('E24').Value = '  +5.29%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.388'
$ws.Range('E25').Value = '  +3.32%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '156.98'
$ws.Range('E26').Value = '  +4.86%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.01'
$ws.Range('E27').Value = '  +11.39%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.914'
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.403'
$ws.Range('E29').Value = '  +0.99%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.313'
$ws.Range('E30').Value = '  +2.40%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08992'
$ws.Range('E31').Value = '  +2.56%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.021'
$ws.Range('E32').Value = '  +1.25%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05191'
$ws.Range('E33').Value = '  +1.98%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7439'
$ws.Range('E34').Value = '  +4.27%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.120'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.722'
$ws.Range('E36').Value = '  +2.05%  '
$ws.Range('E37').Value = '  +2.25%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.678'
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.163'
$ws.Range('E39').Value = '  -0.66%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9421'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4438'
$ws.Range('E41').Value = '  +4.24%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '105.74'
$ws.Range('E42').Value = '  +3.00%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.000'
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.755'
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '7.654'
$ws.Range('E45').Value = '  +2.67%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.1344'
$ws.Range('E46').Value = '  +5.75%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.05845'
$ws.Range('E47').Value = '  +3.14%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.706'
$ws.Range('E48').Value = '  +5.25%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.429'
$ws.Range('E49').Value = '  +6.82%  '
$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.3948'
$ws.Range('E50').Value = '  +4.01%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '33.54'
$ws.Range('E51').Value = '  +2.76%  '
